$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 with the WeaponData entry, following the same
# Name/Path pattern used by the existing rows (row 7 = RealTimePlayerData).
$ws.Range("A8").Value = "WeaponData"
$ws.Range("B8").Value = "WeaponData.xlsx"

# Copy formatting (style) from the previous data row so the new cells
# pick up the same cell style (s="1") as the rest of the table.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
